$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.1358715
$ws.Range("H2").Value = 6.271743000000001
$ws.Range("I2").Value = 0.1851936164128898
$ws.Range("J2").Value = 0.1557817298248122
$ws.Range("O2").Value = 0.03909389944656898
$ws.Range("P2").Value = 0.0575165755591434
$ws.Range("Q2").Value = 0.029343394916
$ws.Range("R2").Value = 0.176060369496
$ws.Range("S2").Value = 0.007239940618191981
$ws.Range("T2").Value = 0.008960031634202873
# Row 3
$ws.Range("G3").Value = 3.1358715
$ws.Range("H3").Value = 6.271743000000001
$ws.Range("I3").Value = 0.1851936164128898
$ws.Range("J3").Value = 0.1557817298248122
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.229998
$ws.Range("N3").Value = 0.459996
$ws.Range("O3").Value = 0.9609061005534311
$ws.Range("P3").Value = 0.9424834244408566
$ws.Range("Q3").Value = 0.7212441732570001
$ws.Range("R3").Value = 2.884976693028
$ws.Range("S3").Value = 0.1779536757946979
$ws.Range("T3").Value = 0.1468216981906093
# Row 4
$ws.Range("I4").Value = 0.360757931601364
$ws.Range("J4").Value = 0.4551951821134977
$ws.Range("O4").Value = 0.03909389944656898
$ws.Range("P4").Value = 0.0575165755591434
$ws.Range("S4").Value = 0.01410343430257593
$ws.Range("T4").Value = 0.02618126808618903
# Row 5
$ws.Range("I5").Value = 0.360757931601364
$ws.Range("J5").Value = 0.4551951821134977
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.229998
$ws.Range("N5").Value = 0.459996
$ws.Range("O5").Value = 0.9609061005534311
$ws.Range("P5").Value = 0.9424834244408566
$ws.Range("Q5").Value = 1.404986635952
$ws.Range("R5").Value = 8.429919815712001
$ws.Range("S5").Value = 0.3466544972987881
$ws.Range("T5").Value = 0.4290139140273087
# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.208691
$ws.Range("H6").Value = 0.626073
$ws.Range("I6").Value = 0.01232456145056403
$ws.Range("J6").Value = 0.01555081815957854
$ws.Range("O6").Value = 0.03909389944656898
$ws.Range("P6").Value = 0.0575165755591434
$ws.Range("Q6").Value = 0.001952791250666667
$ws.Range("R6").Value = 0.017575121256
$ws.Range("S6").Value = 0.0004818151660714103
$ws.Range("T6").Value = 0.0008944298076818987
# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.208691
$ws.Range("H7").Value = 0.626073
$ws.Range("I7").Value = 0.01232456145056403
$ws.Range("J7").Value = 0.01555081815957854
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.229998
$ws.Range("N7").Value = 0.459996
$ws.Range("O7").Value = 0.9609061005534311
$ws.Range("P7").Value = 0.9424834244408566
$ws.Range("Q7").Value = 0.047998512618
$ws.Range("R7").Value = 0.287991075708
$ws.Range("S7").Value = 0.01184274628449262
$ws.Range("T7").Value = 0.01465638835189665
# Row 8
$ws.Range("G8").Value = 7.403123000000001
$ws.Range("H8").Value = 14.806246
$ws.Range("I8").Value = 0.437202583434762
$ws.Range("J8").Value = 0.3677673996035402
$ws.Range("O8").Value = 0.03909389944656898
$ws.Range("P8").Value = 0.0575165755591434
$ws.Range("Q8").Value = 0.06927348961866668
$ws.Range("R8").Value = 0.4156409377120001
$ws.Range("S8").Value = 0.01709195383457877
$ws.Range("T8").Value = 0.0211527214274867
# Row 9
$ws.Range("G9").Value = 7.403123000000001
$ws.Range("H9").Value = 14.806246
$ws.Range("I9").Value = 0.437202583434762
$ws.Range("J9").Value = 0.3677673996035402
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.229998
$ws.Range("N9").Value = 0.459996
$ws.Range("O9").Value = 0.9609061005534311
$ws.Range("P9").Value = 0.9424834244408566
$ws.Range("Q9").Value = 1.702703483754
$ws.Range("R9").Value = 6.810813935016001
$ws.Range("S9").Value = 0.4201106296001833
$ws.Range("T9").Value = 0.3466146781760535
# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.05791866666666667
$ws.Range("H10").Value = 0.173756
$ws.Range("I10").Value = 0.003420474129061952
$ws.Range("J10").Value = 0.004315867255313246
$ws.Range("O10").Value = 0.03909389944656898
$ws.Range("P10").Value = 0.0575165755591434
$ws.Range("Q10").Value = 0.0005419642702222223
$ws.Range("R10").Value = 0.004877678432000001
$ws.Range("S10").Value = 0.0001337196716611385
$ws.Range("T10").Value = 0.0002482339050934571
# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.05791866666666667
$ws.Range("H11").Value = 0.173756
$ws.Range("I11").Value = 0.003420474129061952
$ws.Range("J11").Value = 0.004315867255313246
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.229998
$ws.Range("N11").Value = 0.459996
$ws.Range("O11").Value = 0.9609061005534311
$ws.Range("P11").Value = 0.9424834244408566
$ws.Range("Q11").Value = 0.013321177496
$ws.Range("R11").Value = 0.07992706497600001
$ws.Range("S11").Value = 0.003286754457400814
$ws.Range("T11").Value = 0.004067633350219788
# Row 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.01864033333333333
$ws.Range("H12").Value = 0.055921
$ws.Range("I12").Value = 0.001100832971357958
$ws.Range("J12").Value = 0.001389003043258201
$ws.Range("O12").Value = 0.03909389944656898
$ws.Range("P12").Value = 0.0575165755591434
$ws.Range("Q12").Value = 0.0001744238124444445
$ws.Range("R12").Value = 0.001569814312
$ws.Range("S12").Value = 0.00004303585348973576
$ws.Range("T12").Value = 0.00007989069848944045
# Row 13
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.01864033333333333
$ws.Range("H13").Value = 0.055921
$ws.Range("I13").Value = 0.001100832971357958
$ws.Range("J13").Value = 0.001389003043258201
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.229998
$ws.Range("N13").Value = 0.459996
$ws.Range("O13").Value = 0.9609061005534311
$ws.Range("P13").Value = 0.9424834244408566
$ws.Range("Q13").Value = 0.004287239386
$ws.Range("R13").Value = 0.025723436316
$ws.Range("S13").Value = 0.001057797117868222
$ws.Range("T13").Value = 0.001309112344768761

Write-Host "Applied 150 cell updates"
